$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting (styles) of the last existing data row (row 9)
# onto the new row 10 before putting any values in, so the new cells pick
# up the same cellXfs indices (border style for A10/B10) without Excel
# creating brand-new style entries.
$ws.Range("A9:B9").Copy()
$ws.Range("A10:B10").PasteSpecial()

# A10 holds a dd-mm-yyyy-looking string ("01-10-2025"). Assigning it to
# .Value directly would make Excel auto-convert it into a date serial
# number. Instead, compute it with a TRIM() formula (so Excel is forced to
# treat the result as text) and then convert that formula to its plain
# cached value in-place, which keeps the cell a plain shared-string text
# cell (no formula left behind, no extra number-format style created).
$ws.Range("A10").Formula = '=TRIM("01-10-2025 ")'
$ws.Range("A10").Copy()
$ws.Range("A10").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B10").Value = "The price of gold in India today is ₹11,864 per gram for 24 karat gold, ₹10,875 per gram for 22 karat gold and ₹8,898 per gram for 18 karat gold (also called 999 gold)."

$excel.CutCopyMode = 0
